$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("L_RLC")

# Fix column H formulas: they were copy/pasted from column G (D{row}/2);
# correct them to reference column E (E{row}/2) for both shared-formula blocks.
$ws.Range("H2:H9").Formula = "=E2/2"
$ws.Range("H10:H17").Formula = "=E10/2"

# Apply one-decimal numeric formatting to the sigma columns L and M.
$ws.Range("L2:M17").NumberFormat = "0.0"

# Re-write F1's label so the now-unused duplicate shared string
# ("faseL  (deg)" with a double space) gets dropped from the string table.
$ws.Range("F1").Value = "faseL (deg)"
